$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 / column A value correction: 109189 -> 888888
$ws.Range("A3").Value = 888888

# Move the active selection from G6 to A3
$null = $ws.Range("A3").Select()

# Set an explicit (portrait, letter-size) print setup on the sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
